# Update column C (Fitness) values on Sheet1 for rows 2..162 (Generation 0..160)
# according to the run's recorded convergence values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: starting row, ending row (inclusive), new value for column C
$runs = @(
    @(2, 2, 12953),
    @(3, 14, 10814),
    @(15, 16, 10398),
    @(17, 18, 10228),
    @(19, 31, 10124),
    @(32, 32, 9895),
    @(33, 35, 9834),
    @(36, 45, 9781),
    @(46, 48, 9301),
    @(49, 50, 8961),
    @(51, 53, 8532),
    @(54, 65, 8448),
    @(66, 67, 8237),
    @(68, 78, 8153),
    @(79, 87, 8134),
    @(88, 90, 7995),
    @(91, 136, 7708),
    @(137, 137, 7647),
    @(138, 142, 7622),
    @(143, 162, 7581)
)

foreach ($run in $runs) {
    $startRow = $run[0]
    $endRow = $run[1]
    $value = $run[2]
    $rangeAddr = "C" + $startRow + ":C" + $endRow
    $ws.Range($rangeAddr).Value = $value
}
